$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 5 ("Generating two different games") - Content Placeholder 2
#   Para 1: "We use Rolling Horizon Evolution Algorithm to evolve ..."
#           -> split into 3 runs, replacing "Rolling Horizon Evolution
#              Algorithm " with "Random Mutation Hill Climber "
#   Para 2: "Fitness function used was the one " + "that " + "minimizes..."
#           -> merged back into a single run
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$shape5 = $slide5.Shapes.Item(2)
$tr5 = $shape5.TextFrame.TextRange

$oldPhrase = "Rolling Horizon Evolution Algorithm "
$newPhrase = "Random Mutation Hill Climber "
$full5 = $tr5.Text
$pos = $full5.IndexOf($oldPhrase) + 1
$sub5a = $tr5.Characters($pos, $oldPhrase.Length)
$sub5a.Text = $newPhrase

$full5b = $tr5.Text
$fitnessText = "Fitness function used was the one that minimizes the error based on a perfect value (whether it was 1 life or 6 lives remaining of the winner player at the end of the game)."
$pos2 = $full5b.IndexOf("Fitness function used was the one") + 1
$sub5b = $tr5.Characters($pos2, $fitnessText.Length)
$sub5b.Text = $fitnessText

# ---------------------------------------------------------------------------
# Slides 6 and 7 ("Game A" / "Game B") - Text Placeholder 3
#   Paragraph "Is missile wrap able: " + "False" -> merged into a single run,
#   and the stray trailing endParaRPr mark removed.
# ---------------------------------------------------------------------------
function Merge-MissileWrapParagraph($slideIndex) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text.IndexOf("Is missile wrap able:") -ge 0) {
                $paraCount = $tr.Paragraphs().Count
                for ($j = 1; $j -le $paraCount; $j++) {
                    $para = $tr.Paragraphs($j, 1)
                    if ($para.Text -eq "Is missile wrap able: False") {
                        $prevPara = $tr.Paragraphs($j - 1, 1)
                        $nextPara = $tr.Paragraphs($j + 1, 1)
                        $start = $para.Start
                        $len = $para.Length + $nextPara.Length
                        $span = $tr.Characters($start, $len)
                        $nextText = $nextPara.Text
                        $span.Delete()
                        $ins = $prevPara.InsertAfter("`rIs missile wrap able: False`r" + $nextText)
                        break
                    }
                }
            }
        }
    }
}

Merge-MissileWrapParagraph 6
Merge-MissileWrapParagraph 7
